$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values for the rows
# re-annotated by the SGNN re-run after transcript clean-up.
$updates = @(
    @{ Row = 2; I = "%"; J = "Uninterpretable" }
    @{ Row = 4; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 10; I = "aa"; J = "Agree/Accept" }
    @{ Row = 21; I = "ba"; J = "Appreciation" }
    @{ Row = 23; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 26; I = "aa"; J = "Agree/Accept" }
    @{ Row = 41; I = "sv"; J = "Statement-opinion" }
    @{ Row = 44; I = "aa"; J = "Agree/Accept" }
    @{ Row = 45; I = "aa"; J = "Agree/Accept" }
    @{ Row = 50; I = "aa"; J = "Agree/Accept" }
    @{ Row = 53; I = "sv"; J = "Statement-opinion" }
    @{ Row = 60; I = "aa"; J = "Agree/Accept" }
    @{ Row = 63; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 64; I = "aa"; J = "Agree/Accept" }
    @{ Row = 68; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 71; I = "ba"; J = "Appreciation" }
    @{ Row = 89; I = "sv"; J = "Statement-opinion" }
    @{ Row = 101; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 103; I = "sv"; J = "Statement-opinion" }
    @{ Row = 110; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 113; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 115; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 135; I = "sv"; J = "Statement-opinion" }
    @{ Row = 142; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 149; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 150; I = "sv"; J = "Statement-opinion" }
    @{ Row = 153; I = "sv"; J = "Statement-opinion" }
    @{ Row = 170; I = "sv"; J = "Statement-opinion" }
    @{ Row = 180; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 181; I = "ba"; J = "Appreciation" }
    @{ Row = 197; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 204; I = "%"; J = "Uninterpretable" }
    @{ Row = 211; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 215; I = "sv"; J = "Statement-opinion" }
    @{ Row = 228; I = "aa"; J = "Agree/Accept" }
    @{ Row = 229; I = "ba"; J = "Appreciation" }
    @{ Row = 236; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 241; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 245; I = "sv"; J = "Statement-opinion" }
    @{ Row = 258; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 259; I = "ba"; J = "Appreciation" }
    @{ Row = 260; I = "ba"; J = "Appreciation" }
    @{ Row = 261; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 292; I = "sv"; J = "Statement-opinion" }
    @{ Row = 298; I = "sv"; J = "Statement-opinion" }
    @{ Row = 299; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 302; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 308; I = "qy"; J = "Yes-No-Question" }
    @{ Row = 310; I = "%"; J = "Uninterpretable" }
    @{ Row = 314; I = "ba"; J = "Appreciation" }
    @{ Row = 329; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 347; I = "%"; J = "Uninterpretable" }
    @{ Row = 352; I = "sv"; J = "Statement-opinion" }
    @{ Row = 357; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 362; I = "sv"; J = "Statement-opinion" }
    @{ Row = 373; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 375; I = "sv"; J = "Statement-opinion" }
    @{ Row = 376; I = "sv"; J = "Statement-opinion" }
    @{ Row = 381; I = "sv"; J = "Statement-opinion" }
    @{ Row = 384; I = "sv"; J = "Statement-opinion" }
    @{ Row = 393; I = "qy"; J = "Yes-No-Question" }
    @{ Row = 399; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 407; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 426; I = "sv"; J = "Statement-opinion" }
    @{ Row = 438; I = "qy"; J = "Yes-No-Question" }
    @{ Row = 442; I = "sv"; J = "Statement-opinion" }
    @{ Row = 448; I = "ba"; J = "Appreciation" }
    @{ Row = 451; I = "aa"; J = "Agree/Accept" }
    @{ Row = 452; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 459; I = "sv"; J = "Statement-opinion" }
    @{ Row = 468; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 471; I = "sv"; J = "Statement-opinion" }
    @{ Row = 495; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 500; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
